$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the date form "01.16.19" -> "01.16.20" in the s1cDNADate (column A)
# and s2cDNADate (column D) columns, rows 2 through 32. These cells store the
# date as literal text (shared string), not a real date serial, so we briefly
# force a text number format before writing the value (otherwise a string
# like "01.16.20" gets auto-parsed into a date serial number), then restore
# the cell style back to Normal/General so the saved file keeps the original
# (unstyled) look of these cells.
$ws.Range("A2:A32").NumberFormat = "@"
$ws.Range("A2:A32").Value = "01.16.20"
$ws.Range("A2:A32").Style = "Normal"

$ws.Range("D2:D32").NumberFormat = "@"
$ws.Range("D2:D32").Value = "01.16.20"
$ws.Range("D2:D32").Style = "Normal"

# Restore a selection similar to the authored edit (closest reproducible
# approximation: primary range selected with A2 active).
$ws.Range("A2:A32").Select()
$ws.Range("A2").Activate()
